## Add files via upload
## -----------------------------------------------------------------------
## The EURIBOR6M sheet's first column ("months") is re-labelled as
## "maturities" and every month-count is restated as a "<n>M" label
## (e.g. 0 -> "0M", 1 -> "1M", ..., 360 -> "360M") instead of a bare number.
## -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EURIBOR6M")
$ws.Activate()

# Header: "months" -> "maturities"
$ws.Range("A1").Value = "maturities"

# Row labels: plain month-count numbers -> "<n>M" text labels
$maturities = @(
    "0M", "1M", "2M", "3M", "4M", "5M", "6M", "7M", "8M", "9M",
    "10M", "11M", "12M", "15M", "18M", "21M", "24M", "36M", "48M",
    "60M", "72M", "84M", "96M", "108M", "120M", "144M", "180M",
    "240M", "300M", "360M"
)

for ($i = 0; $i -lt $maturities.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $maturities[$i]
}

# Column A is now text, so widen it to fit the longest label ("maturities")
$ws.Columns.Item(1).ColumnWidth = 10.140625

# Move the active selection as recorded in the saved file
[void]$ws.Range("D7").Select()
